$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")
$ws.Columns("N").Insert()
$ws.Activate()
$ws.Range("T7").Select() | Out-Null
